# Rename the six observation-panel labels (O-1..O-6) to lowercase
# letters (a..f) and nudge/resize their text boxes to the new positions,
# per the authored diff on slide 1 (sldId 263).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$EMU_PER_POINT = 12700

# PowerPoint's Shape.Left/Top/Width/Height are backed by single-precision
# floats, so a naive EMU -> point -> EMU round trip can land 1 EMU off the
# intended value. Nudge in tiny point-sized steps until the stored value
# (read back and re-derived in EMU) matches the requested EMU exactly.
# (Written as four explicit functions rather than one generic/dynamic
# -property one, since dynamic `$obj.$propName = value` assignment isn't
# supported by this PowerShell host.)
function Set-EmuLeft([object]$Shape, [double]$TargetEmu) {
    $pts = $TargetEmu / $EMU_PER_POINT
    $Shape.Left = $pts
    $gotEmu = [Math]::Round($Shape.Left * $EMU_PER_POINT)
    $tries = 0
    while (($gotEmu -ne $TargetEmu) -and ($tries -lt 100)) {
        if ($gotEmu -lt $TargetEmu) { $pts += 0.000005 } else { $pts -= 0.000005 }
        $Shape.Left = $pts
        $gotEmu = [Math]::Round($Shape.Left * $EMU_PER_POINT)
        $tries += 1
    }
}

function Set-EmuTop([object]$Shape, [double]$TargetEmu) {
    $pts = $TargetEmu / $EMU_PER_POINT
    $Shape.Top = $pts
    $gotEmu = [Math]::Round($Shape.Top * $EMU_PER_POINT)
    $tries = 0
    while (($gotEmu -ne $TargetEmu) -and ($tries -lt 100)) {
        if ($gotEmu -lt $TargetEmu) { $pts += 0.000005 } else { $pts -= 0.000005 }
        $Shape.Top = $pts
        $gotEmu = [Math]::Round($Shape.Top * $EMU_PER_POINT)
        $tries += 1
    }
}

function Set-EmuWidth([object]$Shape, [double]$TargetEmu) {
    $pts = $TargetEmu / $EMU_PER_POINT
    $Shape.Width = $pts
    $gotEmu = [Math]::Round($Shape.Width * $EMU_PER_POINT)
    $tries = 0
    while (($gotEmu -ne $TargetEmu) -and ($tries -lt 100)) {
        if ($gotEmu -lt $TargetEmu) { $pts += 0.000005 } else { $pts -= 0.000005 }
        $Shape.Width = $pts
        $gotEmu = [Math]::Round($Shape.Width * $EMU_PER_POINT)
        $tries += 1
    }
}

function Set-EmuHeight([object]$Shape, [double]$TargetEmu) {
    $pts = $TargetEmu / $EMU_PER_POINT
    $Shape.Height = $pts
    $gotEmu = [Math]::Round($Shape.Height * $EMU_PER_POINT)
    $tries = 0
    while (($gotEmu -ne $TargetEmu) -and ($tries -lt 100)) {
        if ($gotEmu -lt $TargetEmu) { $pts += 0.000005 } else { $pts -= 0.000005 }
        $Shape.Height = $pts
        $gotEmu = [Math]::Round($Shape.Height * $EMU_PER_POINT)
        $tries += 1
    }
}

function Set-ShapeGeometry([object]$Shape, [double]$OffX, [double]$OffY, [double]$ExtCx, [double]$ExtCy, [string]$NewText) {
    # Set the text first: these boxes use <a:spAutoFit/>, so changing the
    # text after resizing can make PowerPoint recompute (and override) the
    # height. Doing the text change first avoids that clobbering.
    $Shape.TextFrame.TextRange.Text = $NewText
    Set-EmuLeft $Shape $OffX
    Set-EmuTop $Shape $OffY
    Set-EmuWidth $Shape $ExtCx
    Set-EmuHeight $Shape $ExtCy
}

# Shape index 3 -> id 73 : "O-1" -> "a"
Set-ShapeGeometry $s.Shapes.Item(3) 2239902 657498 252000 276999 "a"

# Shape index 14 -> id 122 : "O-2" -> "b"
Set-ShapeGeometry $s.Shapes.Item(14) 4795525 651843 253079 276999 "b"

# Shape index 15 -> id 123 : "O-3" -> "c"
Set-ShapeGeometry $s.Shapes.Item(15) 2239902 3007190 252000 276999 "c"

# Shape index 16 -> id 124 : "O-4" -> "d"
Set-ShapeGeometry $s.Shapes.Item(16) 4795525 3007189 253079 276999 "d"

# Shape index 17 -> id 125 : "O-5" -> "e"
Set-ShapeGeometry $s.Shapes.Item(17) 2239902 5365880 252000 276999 "e"

# Shape index 18 -> id 126 : "O-6" -> "f"
Set-ShapeGeometry $s.Shapes.Item(18) 4795525 5365880 253079 276999 "f"
